# Add 9 new project rows (UC_1..UC_9) to the UCP_DatasetV1.2 sheet and
# copy the "effort driver" columns for those rows into a new
# "UCC Data Points" worksheet, per the Bayesian-analysis commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------
# 1. New data rows 27-35 on the main sheet
# ---------------------------------------------------------------
# Columns: A Project_No, B Project_Name, C Simple_Actors, D Average_Actors,
#          E Complex_Actors, F UAW(formula), G Simple_UC, H Average_UC,
#          I Complex_UC, J UUCW(formula), K UUCP(formula), L TCF, M ECF,
#          N UCP(formula), O Real_Effort_Person_Hours

$rows = @(
    @{ r=27; name="UC_1"; C=3; D=0; E=0; G=0;  H=2; I=1; L=1.03;                M=2.25;               O=2981.4933332999999 },
    @{ r=28; name="UC_2"; C=3; D=0; E=0; G=2;  H=1; I=0; L=1.0049999999999999;  M=2.4499999999999997; O=58.730999999999995 },
    @{ r=29; name="UC_3"; C=1; D=0; E=0; G=18; H=0; I=0; L=1.0160000000000002;  M=2.25;               O=285.51375000000002 },
    @{ r=30; name="UC_4"; C=1; D=2; E=3; G=2;  H=1; I=0; L=1.0199999999999998;  M=2.0833333333333335; O=104.1 },
    @{ r=31; name="UC_5"; C=2; D=0; E=0; G=0;  H=3; I=1; L=1.0287500000000001;  M=2.25;               O=124.90875 },
    @{ r=32; name="UC_6"; C=2; D=2; E=0; G=2;  H=1; I=0; L=1.02125;             M=2.4;                O=115.46024999999999 },
    @{ r=33; name="UC_7"; C=8; D=0; E=0; G=1;  H=5; I=2; L=1.0306250000000001;  M=2.25;               O=9678.3333332999991 },
    @{ r=34; name="UC_8"; C=2; D=0; E=0; G=22; H=2; I=0; L=1.0156250000000002;  M=2.25;               O=9559.8866667000002 },
    @{ r=35; name="UC_9"; C=1; D=2; E=1; G=3;  H=2; I=0; L=1.0129999999999999;  M=2.3899999999999997; O=174.74625 }
)

$i = 1
foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 1).Value = $i          # A Project_No
    $ws.Cells.Item($r, 2).Value = $row.name   # B Project_Name
    $ws.Cells.Item($r, 3).Value = $row.C      # C Simple_Actors
    $ws.Cells.Item($r, 4).Value = $row.D      # D Average_Actors
    $ws.Cells.Item($r, 5).Value = $row.E      # E Complex_Actors
    $ws.Cells.Item($r, 6).Formula = "=1*C$r+2*D$r+3*E$r"       # F UAW
    $ws.Cells.Item($r, 7).Value = $row.G      # G Simple_UC
    $ws.Cells.Item($r, 8).Value = $row.H      # H Average_UC
    $ws.Cells.Item($r, 9).Value = $row.I      # I Complex_UC
    $ws.Cells.Item($r, 10).Formula = "=G$r*5+H$r*10+I$r*15"    # J UUCW
    $ws.Cells.Item($r, 10).Interior.Color = 65535
    $ws.Cells.Item($r, 11).Formula = "=J$r+F$r"                # K UUCP
    $ws.Cells.Item($r, 11).Interior.Color = 65535
    $ws.Cells.Item($r, 12).Value = $row.L     # L TCF
    $ws.Cells.Item($r, 13).Value = $row.M     # M ECF
    $ws.Cells.Item($r, 14).Formula = "=K$r*L$r*M$r"            # N UCP
    $ws.Cells.Item($r, 14).Interior.Color = 65535
    $ws.Cells.Item($r, 15).Value = $row.O     # O Real_Effort_Person_Hours
    $i = $i + 1
}

# Row 26's K-formula is broken out of the K3:K26 shared-formula group
# (the group now only spans K3:K25, and K27:K35 forms a new group).
$ws.Cells.Item(26, 11).Formula = "=J26+F26"

# Column widths added for the new UC driver columns
$ws.Columns.Item(7).ColumnWidth = 12.7109375
$ws.Columns.Item(8).ColumnWidth = 12.28515625
$ws.Columns.Item(9).ColumnWidth = 11.7109375

# Restore selection/scroll position like the saved file
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("H34").Select()

# ---------------------------------------------------------------
# 2. New "UCC Data Points" worksheet with the effort-driver data
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws)
$ws2.Name = "UCC Data Points"

$ws2.Range("A1").Value = "Proj."
$ws2.Range("B1").Value = "Simple_UC"
$ws2.Range("C1").Value = "Average_UC"
$ws2.Range("D1").Value = "Complex_UC"
$ws2.Range("E1").Value = "EF"
$ws2.Range("F1").Value = "TCF"
$ws2.Range("G1").Value = "Effort"

$data2 = @(
    @{ A=1; B=0;  C=2; D=1; E=1.03;               F=2.25;               G=2981.4933332999999 },
    @{ A=2; B=2;  C=1; D=0; E=1.0049999999999999; F=2.4499999999999997; G=58.730999999999995 },
    @{ A=3; B=18; C=0; D=0; E=1.0160000000000002; F=2.25;               G=285.51375000000002 },
    @{ A=4; B=2;  C=1; D=0; E=1.0199999999999998; F=2.0833333333333335; G=104.1 },
    @{ A=5; B=0;  C=3; D=1; E=1.0287500000000001; F=2.25;               G=124.90875 },
    @{ A=6; B=2;  C=1; D=0; E=1.02125;             F=2.4;               G=115.46024999999999 },
    @{ A=7; B=1;  C=5; D=2; E=1.0306250000000001; F=2.25;               G=9678.3333332999991 },
    @{ A=8; B=22; C=2; D=0; E=1.0156250000000002; F=2.25;               G=9559.8866667000002 },
    @{ A=9; B=3;  C=2; D=0; E=1.0129999999999999; F=2.3899999999999997; G=174.74625 }
)

$r2 = 2
foreach ($d in $data2) {
    $ws2.Cells.Item($r2, 1).Value = $d.A
    $ws2.Cells.Item($r2, 2).Value = $d.B
    $ws2.Cells.Item($r2, 3).Value = $d.C
    $ws2.Cells.Item($r2, 4).Value = $d.D
    $ws2.Cells.Item($r2, 5).Value = $d.E
    $ws2.Cells.Item($r2, 6).Value = $d.F
    $ws2.Cells.Item($r2, 7).Value = $d.G
    $r2 = $r2 + 1
}

$ws2.Columns.Item(2).ColumnWidth = 16
$ws2.Columns.Item(3).ColumnWidth = 15

$ws.Select()
